$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$newChart = $ws.Shapes.AddChart2(-1, 51, 300, 300, 300, 200)
Write-Host $newChart.Chart.ChartType
$newChart.Chart.SetSourceData($ws.Range("C8:D13"))
